$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new column before column C ("HD 511" col 2 shifts to D)
# ------------------------------------------------------------------
$ws.Columns("C").Insert()

# ------------------------------------------------------------------
# 2. Populate new column C ("HD 512") - mirrors column B's data,
#    except row2 (label) and row3 (value)
# ------------------------------------------------------------------
$ws.Range("C2").Value2 = "HD 512"
$ws.Range("C3").Value2 = 3
$ws.Range("C1").Formula = "=C2&"" | ""&C3"

for ($r = 4; $r -le 24; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value2 = $b
}
for ($r = 25; $r -le 29; $r++) {
    $ws.Cells.Item($r, 3).Value2 = "#N/A"
}

# ------------------------------------------------------------------
# 3. Insert three new columns after D (old "HD 511" col) -> E, F, G
#    each a copy of column D, with row 3 changed to 5 / 6 / 8
# ------------------------------------------------------------------
$ws.Columns("E:G").Insert()

$ws.Range("E2").Value2 = "HD 511"
$ws.Range("E3").Value2 = 5
$ws.Range("E1").Formula = "=E2&"" | ""&E3"

$ws.Range("F2").Value2 = "HD 511"
$ws.Range("F3").Value2 = 6
$ws.Range("F1").Formula = "=F2&"" | ""&F3"

$ws.Range("G2").Value2 = "HD 511"
$ws.Range("G3").Value2 = 8
$ws.Range("G1").Formula = "=G2&"" | ""&G3"

for ($r = 4; $r -le 24; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value2 = $d
    $ws.Cells.Item($r, 6).Value2 = $d
    $ws.Cells.Item($r, 7).Value2 = $d
}
for ($r = 25; $r -le 29; $r++) {
    $ws.Cells.Item($r, 5).Value2 = "#N/A"
    $ws.Cells.Item($r, 6).Value2 = "#N/A"
    $ws.Cells.Item($r, 7).Value2 = "#N/A"
}

# ------------------------------------------------------------------
# 4. Column widths
# ------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 11.42578125
$ws.Columns("D").ColumnWidth = 12.42578125
$ws.Columns("E").ColumnWidth = 12.42578125
$ws.Columns("F").ColumnWidth = 9.85546875

# ------------------------------------------------------------------
# 5. Selection
# ------------------------------------------------------------------
$ws.Range("L11").Select()

$excel.Calculate()
